$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revise the value reported for March 2024 (C316)
$ws.Range("C316").Value = 135.56394811247927

# Insert 3 new rows above the trailing "Source" note row (old row 317),
# which pushes that note row down to row 320.
$ws.Range("A317:C319").Insert()

# Fill column A ("year") by copying the existing "2024" cells down, so the
# new rows reuse the same shared-string entry/style as their neighbours
# instead of Excel auto-typing a freshly entered "2024" string as a number.
$ws.Range("A314:A316").Copy()
$ws.Range("A317:A319").PasteSpecial()
$excel.CutCopyMode = $false

# New monthly data rows for 2024: April, May, June
$ws.Range("B317").Value = 4
$ws.Range("C317").Value = 102.66584709712725

$ws.Range("B318").Value = 5
$ws.Range("C318").Value = 144.78526009718755

$ws.Range("B319").Value = 6
$ws.Range("C319").Value = 147.84086669715819
